# Apply updated odds values to Sheet1, matching the target diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2 updates
$ws.Range("G2").Value = 1.62
$ws.Range("H2").Value = 4.75
$ws.Range("I2").Value = 4.5
$ws.Range("L2").Value = 4.75
$ws.Range("N2").Value = 19
$ws.Range("O2").Value = 1.13
$ws.Range("P2").Value = 6
$ws.Range("X2").Value = 10
$ws.Range("Z2").Value = 13
$ws.Range("AA2").Value = 11
$ws.Range("AD2").Value = 9.5
$ws.Range("AQ2").Value = 21
$ws.Range("AR2").Value = 34
$ws.Range("AY2").Value = 23

# Row 4 updates
$ws.Range("M4").Value = 1.04
$ws.Range("N4").Value = 6.6
